# Leads.xlsx: append two more leads (Nhung18, Nhung19) as rows 22-23,
# mirroring the existing Nhung-row pattern, and select the new rows.
# (commit: "thêm dialog comfirm khi remove account")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 22 : lead #21 (Nhung18) ----------------------------------------
$ws.Cells.Item(22, 1).Value2  = 21
$ws.Cells.Item(22, 2).Value2  = "Nhung18"
$ws.Cells.Item(22, 3).Value2  = "Nguyen"
$ws.Cells.Item(22, 4).Formula = '=B22&"@Gmail.com"'
$ws.Cells.Item(22, 5).Value2  = 947948010
$ws.Cells.Item(22, 6).Value2  = "Assigned"
$ws.Cells.Item(22, 7).Value2  = "partner"
$ws.Cells.Item(22, 8).Value2  = "so 10 Dich Vong, Cau Giay, Ha Noi"
$ws.Cells.Item(22, 10).Value2 = "partner business"
$ws.Cells.Item(22, 11).Value2 = "LocDV"
$ws.Cells.Item(22, 13).Formula = '=B22&C22'
$ws.Cells.Item(22, 14).Value2 = 33.78
$ws.Cells.Item(22, 15).Value2 = "buy"

# ---- Row 23 : lead #22 (Nhung19) ----------------------------------------
$ws.Cells.Item(23, 1).Value2  = 22
$ws.Cells.Item(23, 2).Value2  = "Nhung19"
$ws.Cells.Item(23, 3).Value2  = "Nguyen"
$ws.Cells.Item(23, 4).Formula = '=B23&"@Gmail.com"'
$ws.Cells.Item(23, 5).Value2  = 947948011
$ws.Cells.Item(23, 6).Value2  = "Assigned"
$ws.Cells.Item(23, 7).Value2  = "partner"
$ws.Cells.Item(23, 8).Value2  = "so 10 Dich Vong, Cau Giay, Ha Noi"
$ws.Cells.Item(23, 10).Value2 = "partner business"
$ws.Cells.Item(23, 11).Value2 = "LocDV"
$ws.Cells.Item(23, 13).Formula = '=B23&C23'
$ws.Cells.Item(23, 14).Value2 = 34.78
$ws.Cells.Item(23, 15).Value2 = "buy"

# ---- website hyperlinks for the two new rows (column I) -----------------
$ws.Hyperlinks.Add($ws.Range("I22"), "http://www.bkav.com.vn/")
$ws.Hyperlinks.Add($ws.Range("I23"), "http://www.bkav.com.vn/")

# Adding a hyperlink auto-applies the built-in "Hyperlink" cell style;
# the source rows have no such styling, so clear it back off again.
$ws.Range("I22:I23").Style = "Normal"
$hyperlinkStyle = $wb.Styles.Item("Hyperlink")
$hyperlinkStyle.Delete()

# ---- move the selection onto the freshly added rows ----------------------
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A21:XFD23").Select()
